$wb = $excel.ActiveWorkbook

# --- Sheet "rra": add a new data row (row 2) ---
$wsRra = $wb.Worksheets.Item("rra")
$wsRra.Cells.Item(2, 1).Value = "BBMN"
$wsRra.Cells.Item(2, 2).Value = "Jammu & Kashmir"
$wsRra.Cells.Item(2, 3).Value = "DPJ"
$wsRra.Cells.Item(2, 4).Value = "Tamil Nadu"
$wsRra.Cells.Item(2, 5).Value = "RRA"
$wsRra.Cells.Item(2, 6).Value = 1

# --- Sheet "frk_rra": remove the existing data rows (rows 2-4), leaving only the header ---
$wsFrkRra = $wb.Worksheets.Item("frk_rra")
$wsFrkRra.Range("A2:F4").Delete()
